# Regen save_data: update column G (K) values for rows 2-39
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 9
    3  = 5
    4  = 4
    5  = 2
    6  = 2
    7  = 2
    8  = 2
    9  = 3
    10 = 4
    11 = 4
    12 = 5
    13 = 7
    14 = 1
    15 = 3
    16 = 3
    17 = 2
    18 = 6
    19 = 1
    20 = 4
    21 = 3
    22 = 1
    23 = 2
    24 = 2
    25 = 10
    26 = 5
    27 = 2
    28 = 0
    29 = 7
    30 = 5
    31 = 6
    32 = 4
    33 = 4
    34 = 3
    35 = 6
    36 = 6
    37 = 2
    38 = 3
    39 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
